$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item(2)   # female_conc_black_white
$ws3 = $wb.Worksheets.Item(3)   # female_conc_black_sanity
$ws4 = $wb.Worksheets.Item(4)   # female_conc_white_sanity

# --- female_conc_black_white: replace "niece" with "aunt"/"mother"/"grandmother" ---
$ws2.Range("C10").Value = "aunt"
$ws2.Range("C11").Value = "mother"
$ws2.Range("C12").Value = "grandmother"

# --- female_conc_black_sanity: same edit ---
$ws3.Range("C10").Value = "aunt"
$ws3.Range("C11").Value = "mother"
$ws3.Range("C12").Value = "grandmother"

# --- female_conc_white_sanity: same edit, but the whole C2:C12 column loses its
#     inherited cell style (s="1") in the target file, so clear formatting on that
#     range before writing the values back.
$ws4.Range("C2:C12").ClearFormats()
$ws4.Range("C2").Value = "female"
$ws4.Range("C3").Value = "woman"
$ws4.Range("C4").Value = "girl"
$ws4.Range("C5").Value = "sister"
$ws4.Range("C6").Value = "she"
$ws4.Range("C7").Value = "her"
$ws4.Range("C8").Value = "hers"
$ws4.Range("C9").Value = "daughter"
$ws4.Range("C10").Value = "aunt"
$ws4.Range("C11").Value = "mother"
$ws4.Range("C12").Value = "grandmother"

# --- restore view/selection state on each touched sheet, matching the edit session ---
$ws2.Range("C2:C12").Select()
$ws3.Range("C2:C12").Select()
$ws4.Range("C39").Select()
